$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; existing B:E shift right to C:F.
$ws.Range("B1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("B3").Value = "Número da casa"

# Sequential house numbers for the 30 data rows.
for ($i = 0; $i -lt 30; $i++) {
    $ws.Cells.Item(4 + $i, 2).Value = $i + 1
}

# Give the new header cell the same look (bold, bottom border, centered) as
# the other header cells, by copying formatting from an existing header cell.
$ws.Range("C3").Copy()
$ws.Range("B3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F3").PasteSpecial(-4122)   # xlPasteFormats (replaces its old bold/no-border look)

# Give the new data cells the same centered look as the rest of the table.
$ws.Range("C4").Copy()
$ws.Range("B4:B33").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# Match the original author's autofit width for the new column as closely as
# the engine allows.
$ws.Columns.Item(2).ColumnWidth = 13.67

# Mirror the saved selection state from the authored workbook.
$null = $ws.Range("G14").Select()
